$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 25: new entry for 2024-08-23 (serial 45527) - "Créer documents..."
$ws.Range("A25").Value2 = 45527
$ws.Range("B25").Value2 = 4.5
$ws.Range("C25").Value2 = "Créer documents + corriger labo 1 et 2 et révision labo 1 et 2"

# Row 26: new entry for 2024-08-23 (serial 45527) - "Terminer révision labo 2..."
$ws.Range("A26").Value2 = 45527
$ws.Range("B26").Value2 = 2
$ws.Range("C26").Value2 = "Terminer révision labo 2 et ajustement script "

# Match the date formatting + alignment already used by the other date cells in column A
$ws.Range("A25:A26").NumberFormat = "yyyy/mm/dd"
$ws.Range("A25:A26").HorizontalAlignment = -4108

# Leave selection where the author ended up after the edit
$ws.Range("F12").Select() | Out-Null
